# Commit: "Generate Report for Handoff"
#
# The localization-status report was regenerated after a new handoff was
# produced for 4d0d7265-761c-48bd-9dd2-77b05fd541c5.md in the zh-cn locale.
# That pushes the "Latest Handoff Datetime" for that file's zh-cn row from
# 2016-08-17 14:41:37 to 2016-08-17 14:41:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

# Row 5 on the zh-cn sheet corresponds to source file
# 4d0d7265-761c-48bd-9dd2-77b05fd541c5.md (column A). Column H holds the
# "Latest Handoff Datetime" for that row.
$ws.Range("H5").Value = "2016-08-17 14:41:52"
